# Fixed Year Column Issue
# The "Year" column (A2:A6) was missing the "2016-17" academic year and
# instead repeated/misaligned the years below it. Correct the sequence so
# it reads: 2018-19, 2017-18, 2016-17, 2015-16, 2014-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "2016-17"
$ws.Range("A5").Value = "2015-16"
$ws.Range("A6").Value = "2014-15"

# Reflect the author's final selection/cursor position on the sheet.
$ws.Range("A7").Select()
